{"js": "// The document ends with a paragraph containing \"Masivo\" followed by two\n// empty trailing paragraphs (right before the section properties). The\n// commit removes those two trailing empty paragraphs.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Remove the last two paragraphs of the document body (the blank\n// paragraphs that follow the \"Masivo\" paragraph).\nfor (let i = 0; i < 2; i++) {\n  const last = items[items.length - 1 - i];\n  last.delete();\n}\n\nawait context.sync();\n", "ps1": "# The document ends with a paragraph containing \"Masivo\" followed by two\n# empty trailing paragraphs, right before the section properties. This\n# script removes those two trailing empty paragraphs so that the\n# \"Masivo\" paragraph becomes the final paragraph of the document body.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Masivo\" paragraph - the anchor immediately before the two\n# blank paragraphs that need to be removed. Search from the end since it\n# is near the end of the document.\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = $count; $i -ge 1; $i--) {\n    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq \"Masivo\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not find the 'Masivo' paragraph\"\n}\n\n# Everything after it (up to the current last paragraph) should be blank\n# paragraphs that need to disappear.\nif ($count -gt $targetIndex) {\n    # Word never allows deleting the document body's very last paragraph\n    # mark directly (it must always end with one), so first append a\n    # throwaway paragraph. That turns the current trailing blank\n    # paragraph(s) into ordinary, deletable paragraphs.\n    $endRange = $d.Content\n    $endRange.Collapse(0)\n    $endRange.InsertParagraphAfter()\n\n    # Delete the original trailing blank paragraphs, working backwards\n    # from the end, stopping right after the \"Masivo\" paragraph.\n    $current = $d.Paragraphs.Count\n    while ($current -gt $targetIndex + 1) {\n        $d.Paragraphs.Item($current - 1).Range.Delete()\n        $current = $d.Paragraphs.Count\n    }\n\n    # Finally, remove the temporary paragraph appended above - it is no\n    # longer the document's original last paragraph, so it can now be\n    # deleted like any other.\n    $d.Paragraphs.Item($d.Paragraphs.Count).Range.Delete()\n}\n"}
